$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 10-16: new averaged-intensity scheme ordering/values ---
# Row 10: Gaussian-Quadrature
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 1.027139640605709
$ws.Range("D10").Value = 0.9445305089422378
$ws.Range("E10").Value = 1.002723585119218
$ws.Range("F10").Value = 0.9843644865386793
$ws.Range("G10").Value = 1.027139640605709
$ws.Range("H10").Value = 0.9445305089422378
$ws.Range("I10").Value = 1.011351016872909
$ws.Range("J10").Value = 0.9849234737043119
$ws.Range("K10").Value = 1.005664761589806
$ws.Range("L10").Value = 0.9635070909453787
$ws.Range("M10").Value = 1.027139640605709
$ws.Range("N10").Value = 0.9736270470307277
$ws.Range("O10").Value = 0.989689555301461
$ws.Range("P10").Value = 0.9905255705397811

# Row 11: Spiral-90deg-10rot-5space
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 1.06952138450653
$ws.Range("D11").Value = 0.8545695964511686
$ws.Range("E11").Value = 1.037720511284058
$ws.Range("F11").Value = 0.964607202315401
$ws.Range("G11").Value = 1.06952138450653
$ws.Range("H11").Value = 0.8545695964511686
$ws.Range("I11").Value = 1.046213705848453
$ws.Range("J11").Value = 0.973420257073854
$ws.Range("K11").Value = 1.022279546371987
$ws.Range("L11").Value = 0.905381742923465
$ws.Range("M11").Value = 1.06952138450653
$ws.Range("N11").Value = 0.9461450538676134
$ws.Range("O11").Value = 0.9816046736392895
$ws.Range("P11").Value = 0.9842142433468647

# Row 12: Spiral-90deg-15rot-5space
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 1.069098302754471
$ws.Range("D12").Value = 0.8551118185074116
$ws.Range("E12").Value = 1.037605103762144
$ws.Range("F12").Value = 0.9647279220353012
$ws.Range("G12").Value = 1.069098302754471
$ws.Range("H12").Value = 0.8551118185074116
$ws.Range("I12").Value = 1.046021112678542
$ws.Range("J12").Value = 0.9735470773462164
$ws.Range("K12").Value = 1.022117768670891
$ws.Range("L12").Value = 0.9057357111738825
$ws.Range("M12").Value = 1.069098302754471
$ws.Range("N12").Value = 0.9463584611347777
$ws.Range("O12").Value = 0.981635786764832
$ws.Range("P12").Value = 0.9842456021161075

# Row 13: Spiral-90deg-10rot-3space
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 1.069438146469196
$ws.Range("D13").Value = 0.8546784179401092
$ws.Range("E13").Value = 1.037693534019703
$ws.Range("F13").Value = 0.9646337203556069
$ws.Range("G13").Value = 1.069438146469196
$ws.Range("H13").Value = 0.8546784179401092
$ws.Range("I13").Value = 1.046179152124324
$ws.Range("J13").Value = 0.9734574574172539
$ws.Range("K13").Value = 1.022238245074014
$ws.Range("L13").Value = 0.9054483696346373
$ws.Range("M13").Value = 1.069438146469196
$ws.Range("N13").Value = 0.9461859759799061
$ws.Range("O13").Value = 0.9816109546961537
$ws.Range("P13").Value = 0.9842208803793554

# Row 14: NoRotation-tilt60deg
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 1.060048000000001
$ws.Range("D14").Value = 0.8355719999999995
$ws.Range("E14").Value = 1.047087999999999
$ws.Range("F14").Value = 0.9629040000000013
$ws.Range("G14").Value = 1.060048000000001
$ws.Range("H14").Value = 0.8355719999999995
$ws.Range("I14").Value = 1.050015999999998
$ws.Range("J14").Value = 0.9743560000000003
$ws.Range("K14").Value = 1.022519999999999
$ws.Range("L14").Value = 0.9045520000000009
$ws.Range("M14").Value = 1.060048000000001
$ws.Range("N14").Value = 0.9413299999999992
$ws.Range("O14").Value = 0.9764030000000001
$ws.Range("P14").Value = 0.9821319999999999

# Row 15: Rotation-NoTilt
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.1
$ws.Range("D15").Value = 0.77
$ws.Range("E15").Value = 1.06
$ws.Range("F15").Value = 0.95
$ws.Range("G15").Value = 1.1
$ws.Range("H15").Value = 0.77
$ws.Range("I15").Value = 1.07
$ws.Range("J15").Value = 0.96
$ws.Range("K15").Value = 1.04
$ws.Range("L15").Value = 0.87
$ws.Range("M15").Value = 1.1
$ws.Range("N15").Value = 0.915
$ws.Range("O15").Value = 0.97
$ws.Range("P15").Value = 0.9775

# Row 16: Rotation-60detTilt
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.056408986828798
$ws.Range("D16").Value = 0.8637200795648003
$ws.Range("E16").Value = 1.0325399273472
$ws.Range("F16").Value = 0.9685999912959987
$ws.Range("G16").Value = 1.056408986828798
$ws.Range("H16").Value = 0.8637200795648003
$ws.Range("I16").Value = 1.037833666457602
$ws.Range("J16").Value = 0.9740902438912006
$ws.Range("K16").Value = 1.021234507776
$ws.Range("L16").Value = 0.9218742659072036
$ws.Range("M16").Value = 1.056399226879999
$ws.Range("N16").Value = 0.9481300034560001
$ws.Range("O16").Value = 0.9803172462591994
$ws.Range("P16").Value = 0.9845377086336005

# --- Add new rows 17-19 (copy formatting from row 16 template row) ---
# Row 17: HexGrid-90degTilt5degRes
$ws.Range("A16:P16").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9947066764352194
$ws.Range("D17").Value = 0.9942114503031864
$ws.Range("E17").Value = 0.9940420254672179
$ws.Range("F17").Value = 0.9938657476590198
$ws.Range("G17").Value = 0.9947066764352194
$ws.Range("H17").Value = 0.9942114503031864
$ws.Range("I17").Value = 0.9937674776198598
$ws.Range("J17").Value = 0.9945025076966275
$ws.Range("K17").Value = 0.9942840764862741
$ws.Range("L17").Value = 0.9937822702444505
$ws.Range("M17").Value = 0.9947066764352194
$ws.Range("N17").Value = 0.9941267378852021
$ws.Range("O17").Value = 0.9942064749661608
$ws.Range("P17").Value = 0.9941452789889819

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Range("A16:P16").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9898247732663007
$ws.Range("D18").Value = 1.000440016616962
$ws.Range("E18").Value = 0.9919262984649203
$ws.Range("F18").Value = 0.9956326231626454
$ws.Range("G18").Value = 0.9898247732663007
$ws.Range("H18").Value = 1.000440016616962
$ws.Range("I18").Value = 0.9911878322374496
$ws.Range("J18").Value = 0.9963115379682754
$ws.Range("K18").Value = 0.9930930095830671
$ws.Range("L18").Value = 0.9977686341456802
$ws.Range("M18").Value = 0.9898247732663007
$ws.Range("N18").Value = 0.996183157540941
$ws.Range("O18").Value = 0.9944559278777071
$ws.Range("P18").Value = 0.9945230906806626

# Row 19: HexGrid-60degTilt5degRes
$ws.Range("A16:P16").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.9788280942718661
$ws.Range("D19").Value = 1.016001911870168
$ws.Range("E19").Value = 0.9896039602709358
$ws.Range("F19").Value = 0.9990957641612837
$ws.Range("G19").Value = 0.9788280942718661
$ws.Range("H19").Value = 1.016001911870168
$ws.Range("I19").Value = 0.9849644699825476
$ws.Range("J19").Value = 0.9986981690986969
$ws.Range("K19").Value = 0.9892660037337917
$ws.Range("L19").Value = 1.00889057281561
$ws.Range("M19").Value = 0.9788280942718661
$ws.Range("N19").Value = 1.002802936070552
$ws.Range("O19").Value = 0.9958824326435632
$ws.Range("P19").Value = 0.9956686182756125

$excel.CutCopyMode = $false